# Generate Report for Handback
# Update the "generated at" timestamps recorded on the handback status report.
# These are plain text cells (not real datetimes) formatted as yyyy-mm-dd HH:mm:ss.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 04:58:00"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for the first file row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-17 04:57:55"
$wsZhCn.Range("K2").Value = "2016-08-17 04:58:25"

# de-de sheet: Correspond Handback DateTime for the first file row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-17 04:58:32"
